# Duplicate the last sheet ("08-01-22") into a new sheet "08-02-22",
# placed right after it and made the active tab - mirrors what Excel's
# own "Move or Copy... (Create a copy)" does, so the clone keeps the
# source sheet's formatting (header fill/font, phonetic settings, page
# setup, etc.) automatically.

$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$src = $wb.Worksheets.Item($lastIndex)

$src.Copy($null, $src)

$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "08-02-22"

# Keep the header row content explicit (A / B / C), matching the source.
$new.Range("A1").Value = "A"
$new.Range("B1").Value = "B"
$new.Range("C1").Value = "C"

# The newly-copied/active sheet should be the selected tab.
$new.Select()
